$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (F column) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 8
$wsExhibition.Range("F3").Value = 129
$wsExhibition.Range("F4").Value = 684
$wsExhibition.Range("F5").Value = 60

# Sheet "全部类型" (All Types) - update "想去人数" (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8
$wsAll.Range("F4").Value = 129
$wsAll.Range("F5").Value = 684
$wsAll.Range("F6").Value = 60
